$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row has one cell per import-template field/column:
#   ... N=language, O=is_locked_lbl ("锁定"), P=is_enabled_lbl ("启用"/"删除"), Q=rem
# Remove the "is_locked_lbl" and "is_enabled_lbl" columns (O:P) entirely so the
# remaining "rem" column shifts left into O, matching the template's detail
# view without the deleted column.
$ws.Range("O1:P1").EntireColumn.Delete()
